# Updates cryptos list cell values (price + 1h volume %) per upstream data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Cells.Item(2, 4)
$cell.NumberFormat = "@"
$cell.Value = '43.801.15'
$cell.Style = "Normal"

$cell = $ws.Cells.Item(2, 5)
$cell.NumberFormat = "@"
$cell.Value = '  -0.50%  '
$cell.Style = "Normal"

$cell = $ws.Cells.Item(3, 4)
$cell.NumberFormat = "@"
$cell.Value = '2.312.15'
$cell.Style = "Normal"

$cell = $ws.Cells.Item(3, 5)
$cell.NumberFormat = "@"
$cell.Value = '  +3.04%  '
$cell.Style = "Normal"

$cell = $ws.Cells.Item(4, 5)
$cell.NumberFormat = "@"
$cell.Value = '  -0.17%  '
$cell.Style = "Normal"

$cell = $ws.Cells.Item(5, 5)
$cell.NumberFormat = "@"
$cell.Value = '  +1.50%  '
$cell.Style = "Normal"

$cell = $ws.Cells.Item(6, 4)
$cell.NumberFormat = "@"
$cell.Value = '94.23'
$cell.Style = "Normal"

$cell = $ws.Cells.Item(6, 5)
$cell.NumberFormat = "@"
$cell.Value = '  +7.95%  '
$cell.Style = "Normal"

$cell = $ws.Cells.Item(7, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.624'
$cell.Style = "Normal"

$cell = $ws.Cells.Item(7, 5)
$cell.NumberFormat = "@"
$cell.Value = '  -0.14%  '
$cell.Style = "Normal"

$cell = $ws.Cells.Item(8, 5)
$cell.NumberFormat = "@"
$cell.Value = '  -0.11%  '
$cell.Style = "Normal"

$cell = $ws.Cells.Item(9, 5)
$cell.NumberFormat = "@"
$cell.Value = '  +1.51%  '
$cell.Style = "Normal"

$cell = $ws.Cells.Item(10, 4)
$cell.NumberFormat = "@"
$cell.Value = '45.04'
$cell.Style = "Normal"

$cell = $ws.Cells.Item(10, 5)
$cell.NumberFormat = "@"
$cell.Value = '  -1.82%  '
$cell.Style = "Normal"

$cell = $ws.Cells.Item(11, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.0939'
$cell.Style = "Normal"

$cell = $ws.Cells.Item(11, 5)
$cell.NumberFormat = "@"
$cell.Value = '  +1.07%  '
$cell.Style = "Normal"

$cell = $ws.Cells.Item(12, 4)
$cell.NumberFormat = "@"
$cell.Value = '8.14'
$cell.Style = "Normal"

$cell = $ws.Cells.Item(12, 5)
$cell.NumberFormat = "@"
$cell.Value = '  +7.17%  '
$cell.Style = "Normal"

$cell = $ws.Cells.Item(13, 5)
$cell.NumberFormat = "@"
$cell.Value = '  -0.10%  '
$cell.Style = "Normal"

$cell = $ws.Cells.Item(14, 4)
$cell.NumberFormat = "@"
$cell.Value = '2.659.10'
$cell.Style = "Normal"

$cell = $ws.Cells.Item(14, 5)
$cell.NumberFormat = "@"
$cell.Value = '  +2.97%  '
$cell.Style = "Normal"

$cell = $ws.Cells.Item(15, 4)
$cell.NumberFormat = "@"
$cell.Value = '15.32'
$cell.Style = "Normal"

$cell = $ws.Cells.Item(15, 5)
$cell.NumberFormat = "@"
$cell.Value = '  +1.96%  '
$cell.Style = "Normal"

$cell = $ws.Cells.Item(16, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.856'
$cell.Style = "Normal"

$cell = $ws.Cells.Item(16, 5)
$cell.NumberFormat = "@"
$cell.Value = '  +7.10%  '
$cell.Style = "Normal"

$cell = $ws.Cells.Item(17, 4)
$cell.NumberFormat = "@"
$cell.Value = '2.311.63'
$cell.Style = "Normal"

$cell = $ws.Cells.Item(17, 5)
$cell.NumberFormat = "@"
$cell.Value = '  +3.96%  '
$cell.Style = "Normal"

$cell = $ws.Cells.Item(18, 4)
$cell.NumberFormat = "@"
$cell.Value = '43.791.47'
$cell.Style = "Normal"

$cell = $ws.Cells.Item(18, 5)
$cell.NumberFormat = "@"
$cell.Value = '  -0.62%  '
$cell.Style = "Normal"

$cell = $ws.Cells.Item(19, 5)
$cell.NumberFormat = "@"
$cell.Value = '  +1.19%  '
$cell.Style = "Normal"

$cell = $ws.Cells.Item(20, 4)
$cell.NumberFormat = "@"
$cell.Value = '6.29'
$cell.Style = "Normal"

$cell = $ws.Cells.Item(20, 5)
$cell.NumberFormat = "@"
$cell.Value = '  +3.99%  '
$cell.Style = "Normal"

$cell = $ws.Cells.Item(21, 4)
$cell.NumberFormat = "@"
$cell.Value = '71.67'
$cell.Style = "Normal"

$cell = $ws.Cells.Item(21, 5)
$cell.NumberFormat = "@"
$cell.Value = '  +2.00%  '
$cell.Style = "Normal"

$cell = $ws.Cells.Item(22, 4)
$cell.NumberFormat = "@"
$cell.Value = '239.29'
$cell.Style = "Normal"

$cell = $ws.Cells.Item(22, 5)
$cell.NumberFormat = "@"
$cell.Value = '  +2.26%  '
$cell.Style = "Normal"

$cell = $ws.Cells.Item(23, 4)
$cell.NumberFormat = "@"
$cell.Value = '2.29'
$cell.Style = "Normal"

$cell = $ws.Cells.Item(23, 5)
$cell.NumberFormat = "@"
$cell.Value = '  -4.70%  '
$cell.Style = "Normal"

$cell = $ws.Cells.Item(24, 4)
$cell.NumberFormat = "@"
$cell.Value = '9.73'
$cell.Style = "Normal"

$cell = $ws.Cells.Item(24, 5)
$cell.NumberFormat = "@"
$cell.Value = '  +9.00%  '
$cell.Style = "Normal"

$cell = $ws.Cells.Item(25, 5)
$cell.NumberFormat = "@"
$cell.Value = '  +0.06%  '
$cell.Style = "Normal"

$cell = $ws.Cells.Item(26, 4)
$cell.NumberFormat = "@"
$cell.Value = '11.37'
$cell.Style = "Normal"

$cell = $ws.Cells.Item(26, 5)
$cell.NumberFormat = "@"
$cell.Value = '  +3.58%  '
$cell.Style = "Normal"

$cell = $ws.Cells.Item(27, 5)
$cell.NumberFormat = "@"
$cell.Value = '  -0.56%  '
$cell.Style = "Normal"

$cell = $ws.Cells.Item(28, 4)
$cell.NumberFormat = "@"
$cell.Value = '2.33'
$cell.Style = "Normal"

$cell = $ws.Cells.Item(28, 5)
$cell.NumberFormat = "@"
$cell.Value = '  +2.37%  '
$cell.Style = "Normal"

$cell = $ws.Cells.Item(29, 5)
$cell.NumberFormat = "@"
$cell.Value = '  -4.82%  '
$cell.Style = "Normal"

$cell = $ws.Cells.Item(30, 4)
$cell.NumberFormat = "@"
$cell.Value = '39.15'
$cell.Style = "Normal"

$cell = $ws.Cells.Item(30, 5)
$cell.NumberFormat = "@"
$cell.Value = '  -3.29%  '
$cell.Style = "Normal"

$cell = $ws.Cells.Item(31, 4)
$cell.NumberFormat = "@"
$cell.Value = '22.61'
$cell.Style = "Normal"

$cell = $ws.Cells.Item(31, 5)
$cell.NumberFormat = "@"
$cell.Value = '  +8.99%  '
$cell.Style = "Normal"

$cell = $ws.Cells.Item(32, 4)
$cell.NumberFormat = "@"
$cell.Value = '172.42'
$cell.Style = "Normal"

$cell = $ws.Cells.Item(32, 5)
$cell.NumberFormat = "@"
$cell.Value = '  -1.76%  '
$cell.Style = "Normal"

$cell = $ws.Cells.Item(33, 5)
$cell.NumberFormat = "@"
$cell.Value = '  -1.68%  '
$cell.Style = "Normal"

$cell = $ws.Cells.Item(34, 5)
$cell.NumberFormat = "@"
$cell.Value = '  +3.01%  '
$cell.Style = "Normal"

$cell = $ws.Cells.Item(35, 5)
$cell.NumberFormat = "@"
$cell.Value = '  +1.77%  '
$cell.Style = "Normal"

$cell = $ws.Cells.Item(36, 5)
$cell.NumberFormat = "@"
$cell.Value = '  -1.08%  '
$cell.Style = "Normal"

$cell = $ws.Cells.Item(37, 5)
$cell.NumberFormat = "@"
$cell.Value = '  +2.18%  '
$cell.Style = "Normal"

$cell = $ws.Cells.Item(38, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.0358'
$cell.Style = "Normal"

$cell = $ws.Cells.Item(38, 5)
$cell.NumberFormat = "@"
$cell.Value = '  +0.40%  '
$cell.Style = "Normal"

$cell = $ws.Cells.Item(39, 5)
$cell.NumberFormat = "@"
$cell.Value = '  +3.45%  '
$cell.Style = "Normal"

$cell = $ws.Cells.Item(40, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.237'
$cell.Style = "Normal"

$cell = $ws.Cells.Item(40, 5)
$cell.NumberFormat = "@"
$cell.Value = '  +15.74%  '
$cell.Style = "Normal"

$cell = $ws.Cells.Item(41, 5)
$cell.NumberFormat = "@"
$cell.Value = '  +7.15%  '
$cell.Style = "Normal"

$cell = $ws.Cells.Item(42, 4)
$cell.NumberFormat = "@"
$cell.Value = '12.23'
$cell.Style = "Normal"

$cell = $ws.Cells.Item(42, 5)
$cell.NumberFormat = "@"
$cell.Value = '  -4.04%  '
$cell.Style = "Normal"

$cell = $ws.Cells.Item(43, 5)
$cell.NumberFormat = "@"
$cell.Value = '  +17.43%  '
$cell.Style = "Normal"

$cell = $ws.Cells.Item(44, 4)
$cell.NumberFormat = "@"
$cell.Value = '5.46'
$cell.Style = "Normal"

$cell = $ws.Cells.Item(44, 5)
$cell.NumberFormat = "@"
$cell.Value = '  +1.08%  '
$cell.Style = "Normal"

$cell = $ws.Cells.Item(45, 4)
$cell.NumberFormat = "@"
$cell.Value = '61.93'
$cell.Style = "Normal"

$cell = $ws.Cells.Item(45, 5)
$cell.NumberFormat = "@"
$cell.Value = '  -5.10%  '
$cell.Style = "Normal"

$cell = $ws.Cells.Item(46, 4)
$cell.NumberFormat = "@"
$cell.Value = '9.00'
$cell.Style = "Normal"

$cell = $ws.Cells.Item(46, 5)
$cell.NumberFormat = "@"
$cell.Value = '  +7.31%  '
$cell.Style = "Normal"

$cell = $ws.Cells.Item(47, 5)
$cell.NumberFormat = "@"
$cell.Value = '  +3.39%  '
$cell.Style = "Normal"

$cell = $ws.Cells.Item(48, 4)
$cell.NumberFormat = "@"
$cell.Value = '100.53'
$cell.Style = "Normal"

$cell = $ws.Cells.Item(48, 5)
$cell.NumberFormat = "@"
$cell.Value = '  +0.01%  '
$cell.Style = "Normal"

$cell = $ws.Cells.Item(49, 4)
$cell.NumberFormat = "@"
$cell.Value = '1.21'
$cell.Style = "Normal"

$cell = $ws.Cells.Item(49, 5)
$cell.NumberFormat = "@"
$cell.Value = '  +0.03%  '
$cell.Style = "Normal"

$cell = $ws.Cells.Item(50, 4)
$cell.NumberFormat = "@"
$cell.Value = '2.536.36'
$cell.Style = "Normal"

$cell = $ws.Cells.Item(50, 5)
$cell.NumberFormat = "@"
$cell.Value = '  +2.84%  '
$cell.Style = "Normal"

$cell = $ws.Cells.Item(51, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.429'
$cell.Style = "Normal"

$cell = $ws.Cells.Item(51, 5)
$cell.NumberFormat = "@"
$cell.Value = '  -2.81%  '
$cell.Style = "Normal"

